# Reflects Checkpoint 1 Feedback in UserStories and ProjectPlan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Week 4 time-log entry (row 23) plus follow-up note rows (25-29, row 24 left blank)
# Re-use the same date-column style as the existing log rows (copy from A22)
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122) | Out-Null

$ws.Range("A23").Value = 43516
$ws.Range("B23").Value = 1
$ws.Range("D23").Value = "Week 4: Worked on demo videos"

$ws.Range("D25").Value = "Thurs AM worked on demo videos < 1 hour."
$ws.Range("D26").Value = "plus 2 hrs - some challenges reconciling my project w/diffs w/PW's"
$ws.Range("D27").Value = "plus 1 hour finishing demo and working on activity 1 - not sure how to get path for data dump right?"
$ws.Range("D28").Value = "NOT RESOLVED: getting correct path for mysqldump (don't need to do it yet so defering…)"
$ws.Range("D29").Value = "1:30 -x continue week 4 focused on readings (increasingly skimming hibernate tutorial)"

# match column D wrap-text style used by the rest of the log
$ws.Range("D23").WrapText = $true
$ws.Range("D25:D29").WrapText = $true

# Restore the cursor/selection to the new bottom of the log
$ws.Range("D34").Select()
